# Added invivoPKfit outputs to dashboard script
#
# Appends a new "2.3.0" release row to the Table1 benchmark table on
# Sheet1 (invivoPKfit results), then leaves the selection where Excel
# would after typing the row in - on the new Notes cell, R24.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Grow the table by one row - this pushes Table1's ref/autoFilter from
# A1:R23 to A1:R24 and gives us a fresh (formatted) row 24 to fill in.
$tbl = $ws.ListObjects.Item(1)
$tbl.ListRows.Add() | Out-Null

$r = 24

# Match the left-aligned "Normal+left" style used by every other data
# row in the table (style index 1 in the workbook).
$ws.Range("A$r`:R$r").HorizontalAlignment = -4131

$ws.Cells.Item($r, 1).Value  = "2.3.0"                                  # Version
$ws.Cells.Item($r, 2).Value  = 1023                                     # N.steadystate
$ws.Cells.Item($r, 3).Value  = 0.9999                                   # calc_analytic.units
$ws.Cells.Item($r, 4).Value  = 1                                        # solve_pbtk.units
$ws.Cells.Item($r, 5).Value  = 1                                        # calc_mc.units
$ws.Cells.Item($r, 6).Value  = 1.063                                    # RMSLE.Wetmore
$ws.Cells.Item($r, 7).Value  = 352                                      # N.Wetmore
$ws.Cells.Item($r, 8).Value  = 0.2996                                   # RMSLE.noMC
$ws.Cells.Item($r, 9).Value  = 352                                      # N.noMC
$ws.Cells.Item($r, 10).Value = 1.419                                    # RMSLE.InVivoCss
$ws.Cells.Item($r, 11).Value = 86                                       # N.InVivoCss
$ws.Cells.Item($r, 12).Value = 1.047                                    # RMSLE.InVivoAUC
$ws.Cells.Item($r, 13).Value = 86                                       # N.InVivoAUC
$ws.Cells.Item($r, 14).Value = 1.33                                     # RMSLE.InVivoCmax
$ws.Cells.Item($r, 15).Value = 86                                       # N.InVivoCmax
$ws.Cells.Item($r, 16).Value = 0.6344                                   # RMSLE.TissuePC
$ws.Cells.Item($r, 17).Value = 863                                      # N.TissuePC
$ws.Cells.Item($r, 18).Value = "Used Caco-2 to replace Fabs=Fgut=1"     # Notes

# Leave the sheet scrolled/selected on the new row's last cell, as it
# would be right after typing the new Notes entry.
$ws.Activate()
$ws.Range("R$r").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 4
$done = $true
